$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -0.01798179105527164
$ws.Range("E4").Value = 0.005160549134421965
$ws.Range("F4").Value = -0.02386443890417379
$ws.Range("G4").Value = 0.001578179003284269
$ws.Range("H4").Value = -0.0223589975343599
$ws.Range("J4").Value = 0.02126898396478375
$ws.Range("C5").Value = -0.01106302834652113
$ws.Range("E5").Value = -0.009485558395422335
$ws.Range("F5").Value = 0.006785291079105995
$ws.Range("G5").Value = -0.01082770528304635
$ws.Range("H5").Value = 0.002320268732810749
$ws.Range("J5").Value = -0.005967774471864777
$ws.Range("C6").Value = 0.007031711417268455
$ws.Range("E6").Value = 0.0148930026437201
$ws.Range("F6").Value = -0.009373819390042906
$ws.Range("G6").Value = 0.01225754969378855
$ws.Range("H6").Value = -0.000450969426038777
$ws.Range("J6").Value = 0.007038569989143918
$ws.Range("C7").Value = 0.005506116700244667
$ws.Range("E7").Value = 0.01349561823582473
$ws.Range("F7").Value = -0.01745631069777724
$ws.Range("G7").Value = 0.01432388033987983
$ws.Range("H7").Value = 0.009293185619727423
$ws.Range("J7").Value = -0.01252709355526019
$ws.Range("C8").Value = 0.1342772724430909
$ws.Range("E8").Value = -0.004020204544808181
$ws.Range("F8").Value = -0.01448560028589361
$ws.Range("G8").Value = 0.03273622001799333
$ws.Range("H8").Value = 0.9813343973653758
$ws.Range("J8").Value = 0.02624062375787324
$ws.Range("C9").Value = 0.8964337104493483
$ws.Range("E9").Value = -0.001432954713318188
$ws.Range("F9").Value = 0.0169336941335962
$ws.Range("G9").Value = -0.004313662558559949
$ws.Range("H9").Value = -0.02623842402553696
$ws.Range("J9").Value = 0.01603132358852068
$ws.Range("C10").Value = 0.002468269922730797
$ws.Range("E10").Value = 0.003120776380831055
$ws.Range("F10").Value = 0.0152220414021744
$ws.Range("G10").Value = 0.004249022228445377
$ws.Range("H10").Value = 0.01138591111143644
$ws.Range("J10").Value = -0.01661510295137487
$ws.Range("C11").Value = -0.01047466899498676
$ws.Range("E11").Value = 0.006106887604275503
$ws.Range("F11").Value = -0.03207363535907387
$ws.Range("G11").Value = 0.005632428338994987
$ws.Range("H11").Value = 0.005886990667479626
$ws.Range("J11").Value = 0.0102611957578733
$ws.Range("C12").Value = 0.02772284827691392
$ws.Range("E12").Value = 0.002944478133779125
$ws.Range("F12").Value = 0.02263184724506173
$ws.Range("G12").Value = 0.002932572791201252
$ws.Range("H12").Value = 0.02328066179522647
$ws.Range("J12").Value = 0.001637941198684335
$ws.Range("C13").Value = 0.12611549947662
$ws.Range("E13").Value = -0.004301068300042731
$ws.Range("F13").Value = 0.02048141754801347
$ws.Range("G13").Value = -0.004642110784083099
$ws.Range("H13").Value = -0.01053524614940984
$ws.Range("J13").Value = 0.007485079228040229
$ws.Range("C14").Value = -0.2465823403912936
$ws.Range("E14").Value = 0.002288491291539651
$ws.Range("F14").Value = 0.002171666356666391
$ws.Range("G14").Value = 0.005108081380638472
$ws.Range("H14").Value = 0.0198947200597888
$ws.Range("J14").Value = 0.004590071247850292
$ws.Range("C15").Value = -0.008826659681066388
$ws.Range("E15").Value = 0.01296163837446553
$ws.Range("F15").Value = -0.003274600101679226
$ws.Range("G15").Value = 0.01531801402402811
$ws.Range("H15").Value = 0.001885959531438381
$ws.Range("J15").Value = -0.006405788910021091
$ws.Range("C16").Value = -0.004024131712965268
$ws.Range("E16").Value = 0.0008933295717331827
$ws.Range("F16").Value = -0.004246382409524718
$ws.Range("G16").Value = 0.003010070834387564
$ws.Range("H16").Value = 0.0296087301443492
$ws.Range("J16").Value = 0.0007497689453138731
$ws.Range("C17").Value = -0.048725425181017
$ws.Range("E17").Value = -0.002810390512415621
$ws.Range("F17").Value = 0.01409004364871957
$ws.Range("G17").Value = 0.0006606162849821532
$ws.Range("H17").Value = 0.04396874614274984
$ws.Range("J17").Value = -0.009740566160520209
$ws.Range("C18").Value = 0.02792097884483915
$ws.Range("E18").Value = -0.004499464595978583
$ws.Range("F18").Value = -0.006131906291504051
$ws.Range("G18").Value = -0.003576837198683347
$ws.Range("H18").Value = 0.01272460380498415
$ws.Range("J18").Value = -0.008668692563147156
$ws.Range("C19").Value = -0.03826341205853648
$ws.Range("E19").Value = 0.02979383543175341
$ws.Range("F19").Value = -0.009478236714754519
$ws.Range("G19").Value = 0.02986501708496668
$ws.Range("H19").Value = 0.1837101185804047
$ws.Range("J19").Value = -0.007989728808001265
$ws.Range("C20").Value = 0.02801329542453181
$ws.Range("E20").Value = 0.08769243182769726
$ws.Range("F20").Value = -0.05834240862810071
$ws.Range("G20").Value = 0.07082976128448275
$ws.Range("H20").Value = 0.006318580860743234
$ws.Range("J20").Value = -0.1330662151916643
$ws.Range("C21").Value = 0.08081246310449851
$ws.Range("E21").Value = 0.1795254716930189
$ws.Range("F21").Value = -0.07960561263419985
$ws.Range("G21").Value = 0.06182039731663768
$ws.Range("H21").Value = 0.004741184349647373
$ws.Range("J21").Value = 0.8321278226483197
$ws.Range("C22").Value = 0.007189525823581033
$ws.Range("E22").Value = -0.007805346360213854
$ws.Range("F22").Value = 0.01395864660998103
$ws.Range("G22").Value = -0.006763553251093121
$ws.Range("H22").Value = -0.004850024930000997
$ws.Range("J22").Value = 0.002219878705378076
$ws.Range("C23").Value = 0.2218501120740045
$ws.Range("E23").Value = 0.9773427223897089
$ws.Range("F23").Value = -0.5069012008646936
$ws.Range("G23").Value = 0.9866337280407811
$ws.Range("H23").Value = 0.005140930957637237
$ws.Range("J23").Value = -0.4745207522004523
